$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("B2").Value = 3072.051250592643
$ws.Range("C2").Value = 110.4162943388196
$ws.Range("B3").Value = 2719.61491161682
$ws.Range("C3").Value = 94.74019760969304
$ws.Range("B4").Value = 2566.702677368654
$ws.Range("C4").Value = 109.067654123427
$ws.Range("B5").Value = 2744.880363601355
$ws.Range("C5").Value = 75.02382620316575
$ws.Range("B6").Value = 2724.003027639434
$ws.Range("C6").Value = 98.29175387152873
$ws.Range("B7").Value = 2857.737380393959
$ws.Range("C7").Value = 117.9402606970325
$ws.Range("B8").Value = 2674.069619244212
$ws.Range("C8").Value = 108.1670912807062
$ws.Range("B9").Value = 2643.807560029058
$ws.Range("C9").Value = 106.5406700604794
$ws.Range("B10").Value = 2697.469960109578
$ws.Range("C10").Value = 112.4184249422105
$ws.Range("B11").Value = 2278.674148800817
$ws.Range("C11").Value = 117.333143203517
$ws.Range("B12").Value = 1980.497015977507
$ws.Range("C12").Value = 129.2881534066698
$ws.Range("B13").Value = 1858.563128864974
$ws.Range("C13").Value = 116.3759484393393
$ws.Range("B14").Value = 1770.193541444429
$ws.Range("C14").Value = 132.2841217718235
$ws.Range("B15").Value = 2542.044274517867
$ws.Range("C15").Value = 93.89162811493634
$ws.Range("B16").Value = 2645.659933496998
$ws.Range("C16").Value = 99.28500428587859
$ws.Range("B17").Value = 1770.193541444429
$ws.Range("C17").Value = 123.6814502946445
$ws.Range("B18").Value = 1768.506503930825
$ws.Range("C18").Value = 156.1912000440261
$ws.Range("B20").Value = 2154.158497698446
$ws.Range("C20").Value = 110.5301175239406
$ws.Range("B21").Value = 1694.096402980307
$ws.Range("C21").Value = 125.1486714083071
$ws.Range("B22").Value = 1853.789272167813
$ws.Range("C22").Value = 116.9244973969526
$ws.Range("B44").Value = 1507.264949860935
$ws.Range("C44").Value = 300.7760900992207
$ws.Range("B60").Value = 2263.07597718635
$ws.Range("C60").Value = 150.6210631004269
$ws.Range("B193").Value = 4458.022641641485
$ws.Range("C193").Value = 91.25625451668566
